$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New day column (BH): copy the look of the previous day column (BG) ---
# Row 12's player stopped before column BG, so that row must stay untouched;
# copy formats in the two surrounding contiguous blocks instead of the
# whole BG1:BG29 span.
$ws.Range("BG1:BG11").Copy()
$ws.Range("BH1:BH11").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("BG13:BG29").Copy()
$ws.Range("BH13:BH29").PasteSpecial(-4122)  # xlPasteFormats

# Header date for the new day
$ws.Range("BH1").Value = 45938

# Per-player attendance mark for the new day (mirrors BG, with a few
# corrections entered for that day)
$bh = @{
    2  = "P"
    3  = "R"
    4  = "P"
    5  = "P"
    6  = "P"
    7  = "P"
    8  = "B"
    9  = "P"
    10 = "P"
    11 = "P"
    13 = "B"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "B"
    18 = "P"
    19 = "P"
    20 = "P"
    21 = "M"
    22 = "P"
    23 = "B"
    24 = "P"
    25 = "P"
    26 = "P"
    27 = "P"
    28 = "P"
    29 = "P"
}

foreach ($row in $bh.Keys) {
    $ws.Range("BH$row").Value = $bh[$row]
}

# Remove the now-obsolete standalone summary row 30 (was
# "=COUNTIF(BG5:BG29,""P"")"), no longer needed once the per-day column
# totals recompute automatically.
$ws.Rows("30").Delete()

# Move the selection to track the newly added column
$ws.Range("BJ24").Select()
